$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range('A2').Value = 'DTaP '
$ws.Range('A3').Value = 'DTaP '
$ws.Range('A4').Value = 'DTaP '
$ws.Range('A5').Value = 'DTaP-IPV '
$ws.Range('A6').Value = 'DTaP-IPV '
$ws.Range('A7').Value = 'DTaP-IPV '
$ws.Range('A8').Value = 'DTaP-Hep B-IPV '
$ws.Range('A9').Value = 'DTaP-IP-HI '
$ws.Range('A10').Value = 'e-IPV '
$ws.Range('A11').Value = 'Hepatitis A Pediatric '
$ws.Range('A12').Value = 'Hepatitis A Pediatric '
$ws.Range('A13').Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range('A14').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A15').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('B15').Value = 'Recombivax HB'
$ws.Range('A16').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('B16').Value = 'Recombivax HB'
$ws.Range('A17').Value = 'Hib '
$ws.Range('A18').Value = 'Hib '
$ws.Range('A19').Value = 'Hib '
$ws.Range('A20').Value = 'HPV - Human Papillomavirus 9-valent '
$ws.Range('A21').Value = 'MENB - Meningococcal Group B '
$ws.Range('A22').Value = 'MENB - Meningococcal Group B '
$ws.Range('A23').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A24').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A25').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A26').Value = 'MMR/Varicella '
$ws.Range('A27').Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range('A29').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A30').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A31').Value = 'Rotavirus, Live, Oral, Oral '
$ws.Range('A32').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A33').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A34').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A35').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A36').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A37').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A38').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A39').Value = 'Varicella '
$ws = $wb.Worksheets.Item(2)
$ws.Range('A2').Value = 'Hepatitis A Adult '
$ws.Range('A3').Value = 'Hepatitis A Adult '
$ws.Range('A4').Value = 'Hepatitis A Adult '
$ws.Range('A5').Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range('A6').Value = 'Hepatitis B Adult '
$ws.Range('A7').Value = 'Hepatitis B Adult '
$ws.Range('A8').Value = 'Hepatitis B Adult '
$ws.Range('A9').Value = 'HPV-Human Papillomavirus 9 Valent '
$ws.Range('A10').Value = 'Measles, Mumps,  Rubella '
$ws.Range('A11').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A12').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A13').Value = 'MENB - Meningococcal Group B '
$ws.Range('A14').Value = 'MENB - Meningococcal Group B '
$ws.Range('A15').Value = 'Pneumococcal 13-valent '
$ws.Range('A17').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A18').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A19').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A20').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A21').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A22').Value = 'Varicella '
$ws = $wb.Worksheets.Item(3)
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6-35 months)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent Pediatric dose'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluzone Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'Fluarix Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B7').Value = 'FluLaval Quadrivalent'
$ws.Range('A8').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B8').Value = 'FluLaval Quadrivalent'
$ws.Range('A9').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A10').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A11').Value = 'Influenza  (Age 6 -35 months)'
$ws.Range('A12').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('A13').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A14').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B14').Value = 'FluMist Quadrivalent'
$ws = $wb.Worksheets.Item(4)
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluarix Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'FluLaval Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B7').Value = 'FluLaval Quadrivalent'
$ws.Range('A8').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A9').Value = 'Influenza  (Age 4 years and older)'
$ws.Range('A10').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('A11').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B10').Value = 'Afluria Quadrivalent'
$ws.Range('B11').Value = 'Afluria Quadrivalent'
